$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.207.73"

$ws.Range("D3").Value = "3.505.91"
$ws.Range("E3").Value = "  -5.14%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'581.70"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").Value = "'172.64"
$ws.Range("E6").Value = "  -4.47%  "

$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").Value = "3.497.51"
$ws.Range("E8").Value = "  -5.17%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  -6.80%  "

$ws.Range("D11").Value = "'6.69"
$ws.Range("E11").Value = "  +5.34%  "

$ws.Range("D12").Value = "'0.592"
$ws.Range("E12").Value = "  -3.55%  "

$ws.Range("D13").Value = "'46.50"
$ws.Range("E13").Value = "  -6.81%  "

$ws.Range("E14").Value = "  -4.34%  "

$ws.Range("D15").Value = "'673.62"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "4.066.42"
$ws.Range("E16").Value = "  -5.26%  "

$ws.Range("D17").Value = "'8.66"
$ws.Range("E17").Value = "  -4.20%  "

$ws.Range("D18").Value = "69.189.79"
$ws.Range("E18").Value = "  -3.62%  "

$ws.Range("D19").Value = "3.502.70"
$ws.Range("E19").Value = "  -5.21%  "

$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").Value = "'17.37"
$ws.Range("E21").Value = "  -3.99%  "

$ws.Range("D22").Value = "'11.24"
$ws.Range("E22").Value = "  -3.42%  "

$ws.Range("D23").Value = "'0.898"
$ws.Range("E23").Value = "  -5.04%  "

$ws.Range("D24").Value = "'16.07"
$ws.Range("E24").Value = "  -9.71%  "

$ws.Range("D25").Value = "'97.54"
$ws.Range("E25").Value = "  -5.58%  "

$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = "  -4.64%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  -6.87%  "

$ws.Range("D29").Value = "'9.42"
$ws.Range("E29").Value = "  -8.67%  "

$ws.Range("D30").Value = "'32.81"
$ws.Range("E30").Value = "  -8.12%  "

$ws.Range("D31").Value = "'8.66"
$ws.Range("E31").Value = "  -7.02%  "

$ws.Range("D32").Value = "'3.18"
$ws.Range("E32").Value = "  -8.12%  "

$ws.Range("E33").Value = "  -5.28%  "

$ws.Range("D34").Value = "'7.21"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("D35").Value = "'590.15"
$ws.Range("E35").Value = "  +4.15%  "

$ws.Range("D36").Value = "'3.60"
$ws.Range("E36").Value = "  -14.18%  "

$ws.Range("D37").Value = "'10.84"
$ws.Range("E37").Value = "  -4.07%  "

$ws.Range("D38").Value = "'0.103"
$ws.Range("E38").Value = "  -5.45%  "

$ws.Range("D39").Value = "'56.96"
$ws.Range("E39").Value = "  -4.18%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("E41").Value = "  -6.80%  "

$ws.Range("E42").Value = "  -5.50%  "

$ws.Range("D43").Value = "'0.333"
$ws.Range("E43").Value = "  -5.38%  "

$ws.Range("D44").Value = "3.400.70"

$ws.Range("D45").Value = "'33.19"
$ws.Range("E45").Value = "  -6.97%  "

$ws.Range("D46").Value = "0.0₃0704"
$ws.Range("E46").Value = "  -9.36%  "

$ws.Range("D47").Value = "'2.87"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -7.90%  "

$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("D50").Value = "'133.37"
$ws.Range("E50").Value = "  -2.29%  "

$ws.Range("D51").Value = "'5.60"
$ws.Range("E51").Value = "  +13.62%  "
